$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first worker row (EFRAIN IVAN RHENALS HERALDEZ, doc 1143358148);
# everything below shifts up by one row.
$ws.Rows("16:16").Delete()

# Update "Valor Mora" summary value
$ws.Range("E11").Value = 108800

# Update worker / period counts
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

# Reorder periods for MARELIS MERIÑO PEREZ to ascending order (1610, 1611, 1612)
# and refresh the "Valor Mora" amounts for every remaining worker row.
$ws.Range("E16").Value = "1610"
$ws.Range("E17").Value = "1611"
$ws.Range("E18").Value = "1612"

$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
$ws.Range("G19").Value = 1423500
